$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared across Overview!E2/F2/E3/F3, zh-cn!C2/C3, de-de!C2/C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) zh-cn handback report: target file + handback file hyperlinks/names
# ---------------------------------------------------------------------------
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/0f6a8f52d2ba71dd568452aca171e6acd44ae752/e2e/19e67ef5-ce87-4198-97d8-2c3829709996.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/0f6a8f52d2ba71dd568452aca171e6acd44ae752/e2e/51fe2518-3900-478e-8b89-0fabea85b80e.md"

# Rebuild the hyperlink list (existing A2/A3 source links + new I2/I3 target
# links) so relationship ids come out in row order: A2, I2, A3, I3.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $mdUrl1, "", "", "19e67ef5-ce87-4198-97d8-2c3829709996.md")
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", "19e67ef5-ce87-4198-97d8-2c3829709996.md")
$zh.Hyperlinks.Add($zh.Range("A3"), $mdUrl2, "", "", "51fe2518-3900-478e-8b89-0fabea85b80e.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", "51fe2518-3900-478e-8b89-0fabea85b80e.md")

$zh.Range("J2").Value = "19e67ef5-ce87-4198-97d8-2c3829709996.30cdd36964e5361c81f7274a5a2cee28856ae96e.zh-cn.xlf"
$zh.Range("J3").Value = "51fe2518-3900-478e-8b89-0fabea85b80e.5e3d2637957b928b388ecdf0012e3ab778dcceb0.zh-cn.xlf"

# K2/K3 (Latest Handback DateTime) already show the correct handback
# timestamp text once it is refreshed.
$zh.Range("K2").Value = "2016-08-13 08:32:35"
$zh.Range("K3").Value = "2016-08-13 08:32:35"

# ---------------------------------------------------------------------------
# 3) de-de handback report: same shape as zh-cn but with its own timestamp
# ---------------------------------------------------------------------------
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $mdUrl1, "", "", "19e67ef5-ce87-4198-97d8-2c3829709996.md")
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", "19e67ef5-ce87-4198-97d8-2c3829709996.md")
$de.Hyperlinks.Add($de.Range("A3"), $mdUrl2, "", "", "51fe2518-3900-478e-8b89-0fabea85b80e.md")
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", "51fe2518-3900-478e-8b89-0fabea85b80e.md")

$de.Range("J2").Value = "19e67ef5-ce87-4198-97d8-2c3829709996.30cdd36964e5361c81f7274a5a2cee28856ae96e.de-de.xlf"
$de.Range("J3").Value = "51fe2518-3900-478e-8b89-0fabea85b80e.5e3d2637957b928b388ecdf0012e3ab778dcceb0.de-de.xlf"

$de.Range("K2").Value = "2016-08-13 08:32:44"
$de.Range("K3").Value = "2016-08-13 08:32:44"

# ---------------------------------------------------------------------------
# 4) Column widths: widen Status / Latest Target File / Latest Handback File
#    columns so the new, longer content is not clipped.
# ---------------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 29.14
$ov.Columns.Item(6).ColumnWidth = 29.14

$zh.Columns.Item(3).ColumnWidth = 29.14
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

$de.Columns.Item(3).ColumnWidth = 29.14
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17
